$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NBA")

# Fill in new column J (this week's win totals) for rows 2-31
$ws.Range("J2").Value = 10
$ws.Range("J3").Value = 11
$ws.Range("J4").Value = 5
$ws.Range("J5").Value = 11
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 13
$ws.Range("J8").Value = 3
$ws.Range("J9").Value = 7
$ws.Range("J10").Value = 11
$ws.Range("J11").Value = 16
$ws.Range("J12").Value = 13
$ws.Range("J13").Value = 9
$ws.Range("J14").Value = 10
$ws.Range("J15").Value = 16
$ws.Range("J16").Value = 12
$ws.Range("J17").Value = 7
$ws.Range("J18").Value = 9
$ws.Range("J19").Value = 5
$ws.Range("J20").Value = 7
$ws.Range("J21").Value = 10
$ws.Range("J22").Value = 12
$ws.Range("J23").Value = 8
$ws.Range("J24").Value = 4
$ws.Range("J25").Value = 6
$ws.Range("J26").Value = 10
$ws.Range("J27").Value = 7
$ws.Range("J28").Value = 16
$ws.Range("J29").Value = 13
$ws.Range("J30").Value = 11
$ws.Range("J31").Value = 6

# Update the CONCATENATE formulas in column A to include the new $J column
$ws.Range("A2").Formula = '=CONCATENATE("[''",$C2,"'', ","''",$D2,"'' ,",$E2,",",$F2,",",$G2,",",$H2,",",$I2,",",$J2,"],")' 
$ws.Range("A3:A31").Formula = '=CONCATENATE("[''",$C3,"'', ","''",$D3,"'' ,",$E3,",",$F3,",",$G3,",",$H3,",",$I3,",",$J3,"],")' 

# Switch the active sheet to NBA and select A2:A31 (matching the author's final selection)
$ws.Activate()
$ws.Range("A2:A31").Select()
